$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.265.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -8.19%  '
$ws.Range("D3").Value = "'3.177.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -9.92%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = "'511.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -7.81%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = "'170.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -12.97%  '
$ws.Range("D7").Value = "'0.587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -10.40%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = "'3.174.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.89%  '
$ws.Range("D10").Value = "'0.590"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.90%  '
$ws.Range("D11").Value = "'54.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.11%  '
$ws.Range("E12").Value = '  -10.26%  '
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.79%  '
$ws.Range("D14").Value = "'8.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -11.20%  '
$ws.Range("D15").Value = "'3.673.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -10.50%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").Value = "'0.113"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -9.18%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = "'3.163.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -10.59%  '
$ws.Range("D18").Value = "'62.073.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -8.27%  '
$ws.Range("D19").Value = "'16.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.86%  '
$ws.Range("D20").Value = "'10.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -10.19%  '
$ws.Range("D21").Value = "'0.935"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.22%  '
$ws.Range("D22").Value = "'360.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.32%  '
$ws.Range("D23").Value = "'3.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.09%  '
$ws.Range("D24").Value = "'78.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.69%  '
$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").Value = "'10.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.50%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = "'6.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = "'3.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").Value = "'2.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.55%  '
$ws.Range("D29").Value = "'11.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.01%  '
$ws.Range("D30").Value = "'8.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.49%  '
$ws.Range("D31").Value = "'27.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.50%  '
$ws.Range("D32").Value = "'618.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -13.64%  '
$ws.Range("D33").Value = "'6.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.72%  '
$ws.Range("D34").Value = "'10.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.79%  '
$ws.Range("D35").Value = "'0.102"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.72%  '
$ws.Range("D36").Value = "'56.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -11.94%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = "'36.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.03%  '
$ws.Range("D39").Value = "'0.374"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.23%  '
$ws.Range("D40").Value = "'0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("D41").Value = "'0.0₃0662"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("E42").Value = '  -9.91%  '
$ws.Range("D43").Value = "'2.788.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.78%  '
$ws.Range("D44").Value = "'2.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.91%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'2.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.18%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").Value = "'2.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -13.78%  '
$ws.Range("D47").Value = "'0.0381"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.91%  '
$ws.Range("D48").Value = "'2.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.56%  '
$ws.Range("D49").Value = "'2.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.78%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = "'0.121"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.56%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = "'132.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.79%  '
